# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn and de-de sheets, then resize the affected
# status columns to reflect the shorter text (mirrors what Excel's
# column AutoFit would do after the value changes).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns ---
$wsOverview.Range("E2:E4").Value = $newStatus
$wsOverview.Range("F2:F4").Value = $newStatus

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn.Range("C2:C4").Value = $newStatus

# --- de-de sheet: Status column (col C) ---
$wsDeDe.Range("C2:C4").Value = $newStatus

# Resize the status columns to fit the new, shorter text (closest
# achievable snap to the ~13.41 character AutoFit width Excel computes
# for "In Translation").
$wsOverview.Columns("E").ColumnWidth = 12.5
$wsOverview.Columns("F").ColumnWidth = 12.5
$wsZhCn.Columns("C").ColumnWidth = 12.5
$wsDeDe.Columns("C").ColumnWidth = 12.5
